# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Header timestamp text ---
$ws.Range("A1").Value = "Datos actualizados a 7 de Abril de 2020 a las 16:52"

# --- Country rank swap: Irak (row 64) <-> Moldavia (row 65) ---
$ws.Range("A64").Value = "Moldavia"
$ws.Range("A65").Value = "Irak"

# --- Country rank swap: Liechtenstein (row 133) <-> Guatemala (row 134) ---
$ws.Range("A133").Value = "Guatemala"
$ws.Range("A134").Value = "Liechtenstein"

# --- Updated statistics (Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes) ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 369522
$ws.Range("C4").Value = 2518
$ws.Range("D4").Value = 19874
$ws.Range("E4").Value = 338635
$ws.Range("F4").Value = 9015

# Row 15: Paises Bajos
$ws.Range("F15").Value = 1424

# Row 19: Brasil
$ws.Range("B19").Value = 12341
$ws.Range("C19").Value = 158
$ws.Range("E19").Value = 11633
$ws.Range("G19").Value = 17
$ws.Range("H19").Value = 581

# Row 25: Noruega
$ws.Range("E25").Value = 5746
$ws.Range("G25").Value = 12
$ws.Range("H25").Value = 88

# Row 52: Argentina
$ws.Range("E52").Value = 1235
$ws.Range("G52").Value = 2
$ws.Range("H52").Value = 55

# Row 64: now Moldavia (was Irak)
$ws.Range("B64").Value = 1056
$ws.Range("C64").Value = 91
$ws.Range("D64").Value = 40
$ws.Range("E64").Value = 994
$ws.Range("F64").Value = 80
$ws.Range("G64").Value = 3
$ws.Range("H64").Value = 22

# Row 65: now Irak (was Moldavia)
$ws.Range("B65").Value = 1031
$ws.Range("D65").Value = 344
$ws.Range("E65").Value = 623
$ws.Range("F65").Value = 0
$ws.Range("G65").Value = 0
$ws.Range("H65").Value = 64

# Row 71: Barein
$ws.Range("B71").Value = 811
$ws.Range("C71").Value = 55
$ws.Range("E71").Value = 349

# Row 102: Mauricio
$ws.Range("F102").Value = 3

# Row 133: now Guatemala (was Liechtenstein)
$ws.Range("C133").Value = 7
$ws.Range("D133").Value = 17
$ws.Range("E133").Value = 57
$ws.Range("F133").Value = 3
$ws.Range("H133").Value = 3

# Row 134: now Liechtenstein (was Guatemala)
$ws.Range("B134").Value = 77
$ws.Range("C134").Value = 0
$ws.Range("D134").Value = 55
$ws.Range("E134").Value = 21
$ws.Range("F134").Value = 0
$ws.Range("H134").Value = 1
